$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.506.39"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").Value = "3.146.88"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.62"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.15"
$ws.Range("E6").Value = "  -1.89%  "

# Row 8
$ws.Range("D8").Value = "3.142.18"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("E10").Value = "  +0.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.42"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("E12").Value = "  -1.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.42"
$ws.Range("E14").Value = "  +0.02%  "

# Row 15
$ws.Range("D15").Value = "3.662.51"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("E16").Value = "  +2.53%  "

# Row 17
$ws.Range("D17").Value = "64.353.32"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").Value = "3.147.04"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.61"
$ws.Range("E20").Value = "  +1.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.71"
$ws.Range("E21").Value = "  -0.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").Value = "  +1.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("E23").Value = "  -0.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.53"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.39"
$ws.Range("E25").Value = "  -1.66%  "

# Row 26
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("E27").Value = "  -3.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.45"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  +6.43%  "

# Row 30
$ws.Range("E30").Value = "  +1.07%  "

# Row 31
$ws.Range("E31").Value = "  -6.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.07"
$ws.Range("E32").Value = "  +3.56%  "

# Row 33
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  -3.10%  "

# Row 35
$ws.Range("E35").Value = "  -1.72%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0769"
$ws.Range("E37").Value = "  +4.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.46"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.03"
$ws.Range("E39").Value = "  +3.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "447.97"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -0.39%  "

# Row 42
$ws.Range("E42").Value = "  +1.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.24"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("D44").Value = "2.874.02"
$ws.Range("E44").Value = "  +0.92%  "

# Row 45
$ws.Range("E45").Value = "  -1.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  -1.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +2.37%  "

# Row 48
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.21"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("E50").Value = "  -0.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.70"
$ws.Range("E51").Value = "  +0.76%  "
